$d = $word.ActiveDocument

# Change 0
$old0 = "Ativação: 01/01/2018"
$new0 = "Ativação: 01/01/2025"
$d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2) | Out-Null

# Change 1
$old1 = "Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dos Alimentos, abrangendo aspectos sobre as instituições envolvidas, bem como a habilitação dos profissionais nas respectivas áreas de concentração da C&T de Alimentos; segurança alimentar / controle de qualidade; legislação vigente; matérias primas, técnicas de processamento; composição dos alimentos, novos produtos, características do setor alimentício e sua relação com o meio ambiente. Neste contexto, enfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender a demanda por alimentos de qualidade."
$new1 = "Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dosAlimentos, abrangendo aspectos sobre as instituições envolvidas na produção dealimentos, bem como a habilitação dos profissionais nas respectivas áreas deconcentração da C&T de Alimentos; segurança alimentar / controle de qualidade;legislação vigente; produção de matérias primas, técnicas de processamento; boaspráticas de fabricação; composição dos alimentos, desenvolvimento de novos produtos,características do setor alimentício e sua relação com o meio ambiente. Neste contexto,enfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender ademanda por alimentos de qualidade."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Change 3
$old3 = "Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos, segurança alimentar, integridade e características nutritivas dos alimentos); Generalidades sobre o Setor Alimentício Legislação; Matérias-Primas; Técnicas de Processamento: Composição dos Alimentos; Novos produtos"
$new3 = "Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos,segurança alimentar, integridade e características nutritivas dos alimentos);Generalidades sobre o Setor Alimentício Legislação; Produção de Matérias-Primas;Técnicas de Processamento: Composição dos Alimentos; Lançamentos de Novosprodutos"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# Change 4
$old4 = "Food Science and Technolgy: concepts and objectives            (food demand; food safety, nutritional and fisico chemical characteristics of food). Food industries. Legislation; Raw Materials; Food Processing Technics; Food Composition; Development of new foods"
$new4 = "Food Science and Technology: concepts and objectives (food consumption, food safety, integrity and nutritional characteristics of food);General information about the Food Sector Legislation; Production of Raw Materials;Processing Techniques: Food Composition; New Releases products"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# Change 5
$old5 = "Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivas atribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&T Alimentos; embalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos na toxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF); microbiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos de controle - APPCC. Matérias-Primas: importância e características dos segmentos produtores de matérias primas de origem animal, vegetal, microbiana e aditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dos diferentes métodos de processamento dos alimentos abrangendo técnicas de conservação, transformação e melhoria da qualidade; usos e aplicações de aditivos/ingredientes em alimentos. Novos Produtos: fatores que devem ser considerados no lançamento de novos produtos alimentícios. Composição dos Alimentos: composição e significância dos diferentes compostos encontrados nos alimentos focando em suas propriedades e funções"
$new5 = "Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivasatribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&TAlimentos; Eembalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos natoxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF);microbiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos decontrole - APPCC. Matérias-Primas: importância e características dos segmentosprodutores de matérias primas de origem animal, vegetal, microbiana eaditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dosdiferentes métodos de processamento dos alimentos abrangendo técnicas deconservação, transformação e melhoria da qualidade; usos e aplicações deaditivos/ingredientes em alimentos. Novos Produtos: fatores que devem serconsiderados no lançamento de novos produtos alimentícios. Composição dosAlimentos: composição e significância dos diferentes compostos encontrados nosalimentos focando em suas propriedades e funções"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# Change 6
$old6 = "Introduction: concepts of Food Science and Technology;  Major fields of Food Science and professional skills; Institutions involved with Food & Science Technology Development; Food packing; additives and  food ingredients.  Food Safety: aspects of foodborne diseases organic foods; guidelines for food production; food microbiology; industrial hygiene; food quality control. Raw Mateiral for Food Industries: characteristic of different foodstuff. Animal, vegetal and microbial raw materials as foodstuff.  Aditives and food ingredients.Food Processing: objectives and characteristics of different food processing technics. Food preservation, transformations and improvement of food quality methods. Characteristics, use and applications of food additives and ingredients.Novel Food Products:  issues involved in a new food products development.Food Composition: food constituents and its significance. Chemical characteristics, properties and functional aspects of different food compounds. Organic, inorganic and toxic compounds in food composition."
$new6 = "Introduction: concepts of Food Science and Technology; Major fields of Food Scienceand professional skills; Institutions involved with Food & Science TechnologyDevelopment; Food packing; additives and food ingredients.Food Safety: aspects of foodborne diseases organic foods; guidelines for foodproduction; food microbiology; industrial hygiene; food quality control.Raw Material for Food Industries: characteristic of different foodstuff. Animal,vegetal and microbial raw materials as foodstuff. Aditives and food ingredients.Food Processing: objectives and characteristics of different food processing technics.Food preservation, transformations and improvement of food quality methods.Characteristics, use and applications of food additives and ingredients.Novel Food Products: issues involved in a new food products development.Food Composition: food constituents and its significance. Chemical characteristics,properties and functional aspects of different food compounds. Organic, inorganic andtoxic compounds in food composition."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# Change 7
$old7 = "POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.ARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.FRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha, 1981.EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.FENEMA, O.R.,Principles of Food Science: Part I: Food Chemistry. Marcel Dekker, Inc. NY-USA.1975.BENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970."
$new7 = "POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.ARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.FRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha,1981. EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.FENEMA, O.R., Principles of Food Science: Part I: Food Chemistry. Marcel Dekker,Inc. NY-USA.1975.BENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970.TEIXEIRA, L.J.Q; LIMA FILHO, T; SILVA, T; CASSIANO. O. Tecnologia deAlimentos: Processamento Não Térmico Editora: Editora Rubio Ltda. . ISBN:6588340176 ISBN13: 9786588340172/ Páginas: 224 Publicação: Ed.2023ALCARDE, A.R.; D’ARGE, M. R.; SPOLO, M.H.F; Fundamentos de Ciência eTecnologia de Alimentos. E. 2ª - Editora Manole. ISBN: 9788520447147, pags. 480;2019.LAJOLO, F. M & MERCADANTE, A.Z. Química e Bioquímica dos Alimentos -Volume 2 – Ed. 1 - Coleção: Ciência, Tecnologia, Engenharia de Alimentos eNutrição. ISBN-10 ‏ : ‎ 8538808516 ISBN-13 ‏ : ‎ 978-8538808510 . Pag.432. EditoraAtheneu – RJ. 2017."
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# Change 8
$old8 = "LOT2053 -  Microbiologia  (Requisito fraco)"
$new8 = "LOT2053 -  Microbiologia: da Teoria à Prática  (Requisito fraco)"
$d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2) | Out-Null

# Change 2: insert missing English translation text into the empty italic paragraph
$new2 = "The aims of this course are focused on the concepts regarding foodscience and technology and relevant aspects in the food segment; considering issuessuch as food demand; food safety, nutritional and fisico chemical characteristics offood. It is also included aspects regarding Food industries. Legislation; Raw Materials;Food Processing Technics; Food Composition, and Development of new foods"
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim().Length -eq 0 -and $p.Range.Font.Italic -eq -1) {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $targetPara.Range.InsertAfter($new2)
}